$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value (phone number) to the corrected number
$ws.Range("B2").Value = 41998306017

# Remove row 3 entirely (Victor Cals / 41992452510), shifting cells up
$ws.Rows("3:3").Delete()

# Adjust column B width to match the new (narrower) layout
$ws.Columns("B").ColumnWidth = 11.45

# Update the active selection to match the target state
$ws.Range("B3").Select()
